$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab from "BrassA-HW45.xpc" to "BrassA"
$ws.Name = "BrassA"

# Append a new row 16, continuing the pattern of existing row 15:
#   A16 = 14 (bold/bordered/centered style, same as A15)
#   B16 = same text label as B15 ("HexGrid-60degTilt5degRes")
#   C16:P16 = 1
$srcA = $ws.Cells.Item(15, 1)
$dstA = $ws.Cells.Item(16, 1)
$dstA.Value2 = 14
$srcA.Copy()
$dstA.PasteSpecial(-4122)

$ws.Cells.Item(16, 2).Value2 = $ws.Cells.Item(15, 2).Value2

for ($col = 3; $col -le 16; $col++) {
    $ws.Cells.Item(16, $col).Value2 = 1
}

$excel.CutCopyMode = 0
